# Auto-generated Excel COM-interop script
# Applies numeric corrections to the LevePriceNQ/HQ profit calculations
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 31
$ws.Range("H31").Value = 2206.8572
$ws.Range("I31").Value = 74.666664
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 223.999992
$ws.Range("L31").Value = 45000
$ws.Range("M31").Value = 6.000008000000008
$ws.Range("N31").Value = -45460

# Row 51
$ws.Range("H51").Value = 2333
$ws.Range("J51").Value = 2499.5
$ws.Range("L51").Value = 2499.5
$ws.Range("N51").Value = -3467.5

# Row 58
$ws.Range("H58").Value = 1166.8889
$ws.Range("J58").Value = 2512.5
$ws.Range("L58").Value = 7537.5
$ws.Range("N58").Value = -7837.5

# Row 61
$ws.Range("H61").Value = 69420
$ws.Range("I61").Value = 69420
$ws.Range("K61").Value = 208260
$ws.Range("M61").Value = -208088

# Row 112
$ws.Range("H112").Value = 1483.75
$ws.Range("I112").Value = 1665
$ws.Range("J112").Value = 940
$ws.Range("K112").Value = 4995
$ws.Range("L112").Value = 2820
$ws.Range("M112").Value = -3887
$ws.Range("N112").Value = -5036

# Row 113
$ws.Range("H113").Value = 12233.167
$ws.Range("I113").Value = 10680
$ws.Range("J113").Value = 19999
$ws.Range("K113").Value = 10680
$ws.Range("L113").Value = 19999
$ws.Range("M113").Value = -7426
$ws.Range("N113").Value = -26507

# Row 115
$ws.Range("H115").Value = 1030.3334
$ws.Range("I115").Value = 1030.3334
$ws.Range("K115").Value = 3091.0002
$ws.Range("M115").Value = -1524.0002

# Row 131
$ws.Range("H131").Value = 1851.6
$ws.Range("I131").Value = 1814.5
$ws.Range("K131").Value = 5443.5
$ws.Range("M131").Value = -403.5

# Row 137
$ws.Range("H137").Value = 2337.5557
$ws.Range("J137").Value = 3286.375
$ws.Range("L137").Value = 9859.125
$ws.Range("N137").Value = -14959.125

# Row 138
$ws.Range("H138").Value = 1866.5
$ws.Range("I138").Value = 1494.3334
$ws.Range("J138").Value = 4099.5
$ws.Range("K138").Value = 4483.0002
$ws.Range("L138").Value = 12298.5
$ws.Range("M138").Value = 656.9997999999996
$ws.Range("N138").Value = -22578.5


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 55
$ws.Range("H55").Value = 45333
$ws.Range("J55").Value = 45333
$ws.Range("L55").Value = 45333
$ws.Range("N55").Value = -45963

# Row 97
$ws.Range("H97").Value = 1956
$ws.Range("I97").Value = 1565
$ws.Range("J97").Value = 2190.6
$ws.Range("K97").Value = 1565
$ws.Range("L97").Value = 2190.6
$ws.Range("M97").Value = -1069
$ws.Range("N97").Value = -3182.6

# Row 132
$ws.Range("H132").Value = 3133.9092
$ws.Range("I132").Value = 2052.6667
$ws.Range("K132").Value = 6158.000100000001
$ws.Range("M132").Value = -3628.000100000001


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 7061.625
$ws.Range("I86").Value = 3750
$ws.Range("J86").Value = 8165.5
$ws.Range("K86").Value = 3750
$ws.Range("L86").Value = 8165.5
$ws.Range("M86").Value = -2627
$ws.Range("N86").Value = -10411.5

# Row 89
$ws.Range("H89").Value = 7061.625
$ws.Range("I89").Value = 3750
$ws.Range("J89").Value = 8165.5
$ws.Range("K89").Value = 18750
$ws.Range("L89").Value = 40827.5
$ws.Range("M89").Value = -13134
$ws.Range("N89").Value = -52059.5

# Row 94
$ws.Range("H94").Value = 590
$ws.Range("I94").Value = 590
$ws.Range("K94").Value = 590
$ws.Range("M94").Value = -139

# Row 107
$ws.Range("H107").Value = 4940.25
$ws.Range("I107").Value = 1730.5
$ws.Range("K107").Value = 1730.5
$ws.Range("M107").Value = 189.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 85.416664
$ws.Range("I7").Value = 158.75
$ws.Range("J7").Value = 48.75
$ws.Range("K7").Value = 158.75
$ws.Range("L7").Value = 48.75
$ws.Range("M7").Value = -45.75
$ws.Range("N7").Value = -274.75

# Row 47
$ws.Range("H47").Value = 24999.5
$ws.Range("I47").Value = 24999.5
$ws.Range("K47").Value = 24999.5
$ws.Range("M47").Value = -24433.5

# Row 62
$ws.Range("H62").Value = 2502.5
$ws.Range("I62").Value = 2005
$ws.Range("K62").Value = 2005
$ws.Range("M62").Value = -1381

# Row 65
$ws.Range("H65").Value = 2502.5
$ws.Range("I65").Value = 2005
$ws.Range("K65").Value = 10025
$ws.Range("M65").Value = -6905

# Row 122
$ws.Range("H122").Value = 1090.5714
$ws.Range("I122").Value = 941.1667
$ws.Range("J122").Value = 1987
$ws.Range("K122").Value = 2823.5001
$ws.Range("L122").Value = 5961
$ws.Range("M122").Value = -373.5001000000002
$ws.Range("N122").Value = -10861

# Row 132
$ws.Range("H132").Value = 4682.0605
$ws.Range("I132").Value = 3787.9473
$ws.Range("K132").Value = 11363.8419
$ws.Range("M132").Value = -8833.841899999999

# Row 134
$ws.Range("H134").Value = 2214.9
$ws.Range("I134").Value = 2214.9
$ws.Range("K134").Value = 6644.700000000001
$ws.Range("M134").Value = -4109.700000000001


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 37928024
$ws.Range("I4").Value = 56889256
$ws.Range("J4").Value = 5555.4443
$ws.Range("K4").Value = 170667768
$ws.Range("L4").Value = 16666.3329
$ws.Range("M4").Value = -170667656
$ws.Range("N4").Value = -16890.3329

# Row 25
$ws.Range("H25").Value = 37
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 52.5
$ws.Range("K25").Value = 18
$ws.Range("L25").Value = 157.5
$ws.Range("M25").Value = 151
$ws.Range("N25").Value = -495.5

# Row 30
$ws.Range("H30").Value = 37
$ws.Range("I30").Value = 6
$ws.Range("J30").Value = 52.5
$ws.Range("K30").Value = 18
$ws.Range("L30").Value = 157.5
$ws.Range("M30").Value = 84
$ws.Range("N30").Value = -361.5

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Row 122
$ws.Range("H122").Value = 1715.6666
$ws.Range("J122").Value = 1624.5
$ws.Range("L122").Value = 14620.5
$ws.Range("N122").Value = -19520.5

# Row 128
$ws.Range("H128").Value = 629999
$ws.Range("I128").Value = 629999
$ws.Range("K128").Value = 1889997
$ws.Range("M128").Value = -1885017


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 1716.6666
$ws.Range("I102").Value = 1716.6666
$ws.Range("K102").Value = 1716.6666
$ws.Range("M102").Value = -94.66660000000002

# Row 122
$ws.Range("H122").Value = 627131.75
$ws.Range("I122").Value = 716507.7
$ws.Range("K122").Value = 2149523.1
$ws.Range("M122").Value = -2147073.1


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 7549.769
$ws.Range("I7").Value = 5449.5
$ws.Range("J7").Value = 8483.223
$ws.Range("K7").Value = 5449.5
$ws.Range("L7").Value = 8483.223
$ws.Range("M7").Value = -5337.5
$ws.Range("N7").Value = -8707.223

# Row 22
$ws.Range("H22").Value = 649.6667
$ws.Range("I22").Value = 749.5
$ws.Range("J22").Value = 599.75
$ws.Range("K22").Value = 749.5
$ws.Range("L22").Value = 599.75
$ws.Range("M22").Value = -454.5
$ws.Range("N22").Value = -1189.75

# Row 27
$ws.Range("H27").Value = 649.6667
$ws.Range("I27").Value = 749.5
$ws.Range("J27").Value = 599.75
$ws.Range("K27").Value = 749.5
$ws.Range("L27").Value = 599.75
$ws.Range("M27").Value = -642.5
$ws.Range("N27").Value = -813.75

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 122
$ws.Range("H122").Value = 5929
$ws.Range("I122").Value = 6083.8335
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 18251.5005
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -15801.5005
$ws.Range("N122").Value = -19900

# Row 126
$ws.Range("H126").Value = 7549.769
$ws.Range("I126").Value = 5449.5
$ws.Range("J126").Value = 8483.223
$ws.Range("K126").Value = 16348.5
$ws.Range("L126").Value = 25449.669
$ws.Range("M126").Value = -13878.5
$ws.Range("N126").Value = -30389.669

# Row 132
$ws.Range("H132").Value = 2862.6
$ws.Range("I132").Value = 2961.75
$ws.Range("J132").Value = 2466
$ws.Range("K132").Value = 8885.25
$ws.Range("L132").Value = 7398
$ws.Range("M132").Value = -6355.25
$ws.Range("N132").Value = -12458


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 14
$ws.Range("H14").Value = 812.3333
$ws.Range("I14").Value = 921.6667
$ws.Range("K14").Value = 921.6667
$ws.Range("M14").Value = -753.6667

# Row 56
$ws.Range("H56").Value = 31656.5
$ws.Range("J56").Value = 31656.5
$ws.Range("L56").Value = 31656.5
$ws.Range("N56").Value = -33084.5

# Row 113
$ws.Range("H113").Value = 941.1818
$ws.Range("I113").Value = 1006.75
$ws.Range("K113").Value = 3020.25
$ws.Range("M113").Value = -850.25

